# Daily attendance processing - 2026-01-24 15:59:17
# Normalize the "Recorded By" column (G): move the trailing contributor
# (the most recently-appended recorder, e.g. "System") to the front of the
# comma-separated list so the most recent recorder is listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -ge 2) {
        $lastPart = $parts[$parts.Count - 1]
        $rest = $parts[0..($parts.Count - 2)]
        $newParts = @($lastPart) + $rest
        $newVal = $newParts -join ", "
        $cell.Value = $newVal
    }
}
